$wb = $excel.ActiveWorkbook

# ---- Timetable ----
$ws = $wb.Worksheets.Item("Timetable")
$ws.Range("E3").Value = "CS307 [C304]"
$ws.Range("C8").Value = "DS302 (Lab) [L106]"
$ws.Range("E8").Value = "CS307 (Lab) [L207]"
$ws.Range("B9").Value = "DS303 (Tutorial) [C101]"
$ws.Range("C9").Value = "DS302 (Lab) [L106]"
$ws.Range("E9").Value = "CS307 (Lab) [L207]"

# ---- Verification ----
$ws = $wb.Worksheets.Item("Verification")
$ws.Range("I2").Value = "C004, L106, C101"
$ws.Range("I3").Value = "C004, L402, C101"
$ws.Range("C7").Value = "Utkarsh Mahadeo Khaire, Siddharth R, Deepak K T"
$ws.Range("D7").Value = "3-0-2-0-4"
$ws.Range("I7").Value = "C004, L207, C304"

# ---- Room_Allocation ----
$ws = $wb.Worksheets.Item("Room_Allocation")
$ws.Range("H2").Value = "MINOR: Cybersecurity, MINOR: Generative Ai, MINOR: Design..."
$ws.Range("H3").Value = "DS302, CS307, DS303"
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = "DS302 (Tutorial), DS303 (Tutorial)"
$ws.Range("I4").Value = 0.4
$ws.Range("A5").Value = "C304"
$ws.Range("H5").Value = "CS307"
$ws.Range("H6").Value = "DS302 (Lab)"
$ws.Range("H7").Value = "CS307 (Lab)"
$ws.Range("E8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = "DS303"
$ws.Range("I8").Value = 0.2

# ---- Classroom_Allocation ----
$ws = $wb.Worksheets.Item("Classroom_Allocation")
$ws.Range("I6").Value = "Projector"
$ws.Range("M6").Value = "C101"
$ws.Range("M11").Value = "L106"
$ws.Range("M12").Value = "L106"
$ws.Range("H18").Value = "96"
$ws.Range("M18").Value = "C304"
$ws.Range("M19").Value = "L207"
$ws.Range("M20").Value = "L207"
$ws.Range("G25").Value = "classroom"
$ws.Range("H25").Value = "96"
$ws.Range("M25").Value = "C104"
$ws.Range("I26").Value = "TV"
$ws.Range("M26").Value = "C203"
$ws.Range("I27").Value = "TV"
$ws.Range("M27").Value = "C204"
$ws.Range("I28").Value = "TV"
$ws.Range("M28").Value = "C205"
$ws.Range("G30").Value = "classroom"
$ws.Range("H30").Value = "96"
$ws.Range("I30").Value = "Projector"
$ws.Range("M30").Value = "C102"
$ws.Range("G31").Value = "classroom"
$ws.Range("I31").Value = "Projector"
$ws.Range("M31").Value = "C104"
$ws.Range("G32").Value = "classroom"
$ws.Range("H32").Value = "96"
$ws.Range("M32").Value = "C202"
$ws.Range("I33").Value = "TV"
$ws.Range("M33").Value = "C203"
$ws.Range("G34").Value = "classroom"
$ws.Range("H34").Value = "96"
$ws.Range("I34").Value = "Projector"
$ws.Range("M34").Value = "C102"
$ws.Range("G35").Value = "classroom"
$ws.Range("I35").Value = "Projector"
$ws.Range("M35").Value = "C104"
$ws.Range("G36").Value = "classroom"
$ws.Range("H36").Value = "96"
$ws.Range("M36").Value = "C202"
$ws.Range("I37").Value = "TV"
$ws.Range("M37").Value = "C203"
$ws.Range("G39").Value = "classroom"
$ws.Range("H39").Value = "96"
$ws.Range("M39").Value = "C104"
$ws.Range("I40").Value = "TV"
$ws.Range("M40").Value = "C203"
$ws.Range("I41").Value = "TV"
$ws.Range("M41").Value = "C204"
$ws.Range("I42").Value = "TV"
$ws.Range("M42").Value = "C205"
$ws.Range("G44").Value = "classroom"
$ws.Range("H44").Value = "96"
$ws.Range("I44").Value = "Projector"
$ws.Range("M44").Value = "C102"
$ws.Range("G45").Value = "classroom"
$ws.Range("I45").Value = "Projector"
$ws.Range("M45").Value = "C104"
$ws.Range("G46").Value = "classroom"
$ws.Range("H46").Value = "96"
$ws.Range("M46").Value = "C202"
$ws.Range("I47").Value = "TV"
$ws.Range("M47").Value = "C203"
$ws.Range("G49").Value = "classroom"
$ws.Range("H49").Value = "96"
$ws.Range("M49").Value = "C104"
$ws.Range("I50").Value = "TV"
$ws.Range("M50").Value = "C203"
$ws.Range("I51").Value = "TV"
$ws.Range("M51").Value = "C204"
$ws.Range("I52").Value = "TV"
$ws.Range("M52").Value = "C205"

# ---- LTPSC_Compliance ----
$ws = $wb.Worksheets.Item("LTPSC_Compliance")
$ws.Range("C2").Value = "3-0-2-0-4"
$ws.Range("D2").Value = "3/0/2"
$ws.Range("H2").Value = "[OK]"
$ws.Range("I2").Value = "[FAIL]"

# ---- Executive_Summary ----
$ws = $wb.Worksheets.Item("Executive_Summary")
$ws.Range("C3").Value = "2026-01-26 12:46"

# ---- Course_Summary ----
$ws = $wb.Worksheets.Item("Course_Summary")
$ws.Rows.Item(15).Delete()

# ---- Classroom_Utilization ----
$ws = $wb.Worksheets.Item("Classroom_Utilization")
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 0.4
$ws.Range("G6").Value = 5
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("D25").Value = 1.5
$ws.Range("E25").Value = 0.3
$ws.Range("G25").Value = 3.75
$ws.Range("D31").Value = 1.5
$ws.Range("E31").Value = 0.3
$ws.Range("G31").Value = 3.75

# ---- Section_A ----
$ws = $wb.Worksheets.Item("Section_A")
$ws.Range("E3").Value = "CS307 [C304]"
$ws.Range("C8").Value = "DS302 (Lab) [L106]"
$ws.Range("E8").Value = "CS307 (Lab) [L207]"
$ws.Range("B9").Value = "DS303 (Tutorial) [C101]"
$ws.Range("C9").Value = "DS302 (Lab) [L106]"
$ws.Range("E9").Value = "CS307 (Lab) [L207]"

# ---- Basket_Course_Allocations ----
$ws = $wb.Worksheets.Item("Basket_Course_Allocations")
$ws.Range("C5").Value = "C102"
$ws.Range("C6").Value = "C104"
$ws.Range("C7").Value = "C202"
$ws.Range("C8").Value = "C203"
$ws.Range("C10").Value = "C104"
$ws.Range("C11").Value = "C204"
$ws.Range("C12").Value = "C102"
$ws.Range("C13").Value = "C104"
$ws.Range("C14").Value = "C202"
$ws.Range("C15").Value = "C203"
$ws.Range("C17").Value = "C104"
$ws.Range("C18").Value = "C203"
$ws.Range("C19").Value = "C204"
$ws.Range("C20").Value = "C205"
